$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds the numeric ruler 0..21 in columns A..V.
# The edit shifts that ruler one column to the right (A23 becomes blank,
# B23 gets the old A23 value, ..., V23 gets the old U23 value), and the
# old V23 value (21) is dropped. Column W23 (and beyond) is untouched.
for ($col = 22; $col -ge 2; $col--) {
    $srcCell = $ws.Cells.Item(23, $col - 1)
    $dstCell = $ws.Cells.Item(23, $col)
    $dstCell.Value = $srcCell.Value2
}
$ws.Cells.Item(23, 1).Value = $null

# Update the sheet's active selection to U17 (previously V23).
$ws.Range("U17").Select()
